# Update ZNH Yearly Financials worksheet with refreshed data pull
# (values sourced from an updated financial data feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZNH")

$ws.Range("D8").Value = 18967700
$ws.Range("E8").Value = 17064300
$ws.Range("F8").Value = 16570300
$ws.Range("G8").Value = 16115000
$ws.Range("H8").Value = 14625400
$ws.Range("I8").Value = 14768900
$ws.Range("J8").Value = 13415500
$ws.Range("D9").Value = 14714000
$ws.Range("E9").Value = 12622000
$ws.Range("F9").Value = 12103700
$ws.Range("G9").Value = 13239400
$ws.Range("H9").Value = 12115000
$ws.Range("I9").Value = 11387900
$ws.Range("J9").Value = 10123300
$ws.Range("D10").Value = 4253700
$ws.Range("E10").Value = 4442400
$ws.Range("F10").Value = 4466500
$ws.Range("G10").Value = 2875600
$ws.Range("H10").Value = 2510400
$ws.Range("I10").Value = 3380900
$ws.Range("J10").Value = 3292200
$ws.Range("D14").Value = 31900
$ws.Range("F14").Value = 13400
$ws.Range("G14").Value = 28000
$ws.Range("H14").Value = 79500
$ws.Range("J14").Value = 86700
$ws.Range("I15").Value = 1226500
$ws.Range("J15").Value = 1141100
$ws.Range("D17").Value = 17592700
$ws.Range("E17").Value = 15179200
$ws.Range("F17").Value = 14575900
$ws.Range("G17").Value = 15406400
$ws.Range("H17").Value = 14401300
$ws.Range("I17").Value = 14012100
$ws.Range("J17").Value = 12769500
$ws.Range("D18").Value = 1375000
$ws.Range("E18").Value = 1885100
$ws.Range("F18").Value = 1994300
$ws.Range("G18").Value = 708500
$ws.Range("H18").Value = 224100
$ws.Range("I18").Value = 756700
$ws.Range("J18").Value = 646000
$ws.Range("D20").Value = 347000
$ws.Range("E20").Value = -389300
$ws.Range("F20").Value = -761900
$ws.Range("G20").Value = 71400
$ws.Range("H20").Value = 537200
$ws.Range("I20").Value = 150600
$ws.Range("J20").Value = 540800
$ws.Range("D21").Value = 3702700
$ws.Range("E21").Value = 3391400
$ws.Range("F21").Value = 3015400
$ws.Range("G21").Value = 2412900
$ws.Range("H21").Value = 2172700
$ws.Range("I21").Value = 2147000
$ws.Range("J21").Value = 2340800
$ws.Range("D22").Value = 405000
$ws.Range("E22").Value = 358900
$ws.Range("F22").Value = 324400
$ws.Range("G22").Value = 324900
$ws.Range("H22").Value = 244300
$ws.Range("I22").Value = 204200
$ws.Range("J22").Value = 158400
$ws.Range("D23").Value = 1317000
$ws.Range("E23").Value = 1137000
$ws.Range("F23").Value = 908000
$ws.Range("G23").Value = 455000
$ws.Range("H23").Value = 517100
$ws.Range("I23").Value = 703200
$ws.Range("J23").Value = 1028500
$ws.Range("D24").Value = 293300
$ws.Range("E24").Value = 261600
$ws.Range("F24").Value = 192900
$ws.Range("G24").Value = 99100
$ws.Range("H24").Value = 108900
$ws.Range("I24").Value = 141600
$ws.Range("J24").Value = 124700
$ws.Range("D26").Value = 1023700
$ws.Range("E26").Value = 875300
$ws.Range("F26").Value = 715000
$ws.Range("G26").Value = 355900
$ws.Range("H26").Value = 408100
$ws.Range("I26").Value = 561600
$ws.Range("J26").Value = 903800
$ws.Range("D27").Value = 884700
$ws.Range("E27").Value = 748600
$ws.Range("F27").Value = 554500
$ws.Range("G27").Value = 263700
$ws.Range("H27").Value = 294700
$ws.Range("I27").Value = 388700
$ws.Range("J27").Value = 758400
$ws.Range("D32").Value = -347000
$ws.Range("E32").Value = 389300
$ws.Range("F32").Value = 761900
$ws.Range("G32").Value = -71400
$ws.Range("H32").Value = -537200
$ws.Range("I32").Value = -150600
$ws.Range("J32").Value = -540800
$ws.Range("D33").Value = 884700
$ws.Range("E33").Value = 748600
$ws.Range("F33").Value = 554500
$ws.Range("G33").Value = 263700
$ws.Range("H33").Value = 294700
$ws.Range("I33").Value = 388700
$ws.Range("J33").Value = 758400
$ws.Range("D35").Value = 884700
$ws.Range("E35").Value = 748600
$ws.Range("F35").Value = 554500
$ws.Range("G35").Value = 263700
$ws.Range("H35").Value = 294700
$ws.Range("I35").Value = 388700
$ws.Range("J35").Value = 758400
$ws.Range("D41").Value = 1013000
$ws.Range("E41").Value = 612300
$ws.Range("F41").Value = 662200
$ws.Range("G41").Value = 2287600
$ws.Range("H41").Value = 2590300
$ws.Range("I41").Value = 1496300
$ws.Range("J41").Value = 705700
$ws.Range("E42").Value = 3900
$ws.Range("F42").Value = 14500
$ws.Range("H42").Value = 1024600
$ws.Range("I42").Value = 1496300
$ws.Range("J42").Value = 758100
$ws.Range("D43").Value = 1184800
$ws.Range("E43").Value = 960800
$ws.Range("F43").Value = 984400
$ws.Range("G43").Value = 1340600
$ws.Range("H43").Value = 992600
$ws.Range("I43").Value = 629100
$ws.Range("J43").Value = 1058900
$ws.Range("D44").Value = 240700
$ws.Range("E44").Value = 235700
$ws.Range("F44").Value = 238300
$ws.Range("G44").Value = 246500
$ws.Range("H44").Value = 488900
$ws.Range("I44").Value = 482300
$ws.Range("J44").Value = 240100
$ws.Range("D45").Value = 215600
$ws.Range("E45").Value = 230000
$ws.Range("F45").Value = 195000
$ws.Range("G45").Value = 212700
$ws.Range("H45").Value = 381400
$ws.Range("I45").Value = 149200
$ws.Range("J45").Value = 129000
$ws.Range("D46").Value = 2654200
$ws.Range("E46").Value = 2042700
$ws.Range("F46").Value = 2094500
$ws.Range("G46").Value = 4087400
$ws.Range("H46").Value = 3052900
$ws.Range("I46").Value = 2491400
$ws.Range("J46").Value = 2891800
$ws.Range("D47").Value = 803300
$ws.Range("E47").Value = 807200
$ws.Range("F47").Value = 689800
$ws.Range("G47").Value = 610300
$ws.Range("H47").Value = 692300
$ws.Range("I47").Value = 450700
$ws.Range("J47").Value = 352900
$ws.Range("D48").Value = 28099900
$ws.Range("E48").Value = 26086800
$ws.Range("F48").Value = 24104900
$ws.Range("G48").Value = 22825500
$ws.Range("H48").Value = 38118800
$ws.Range("I48").Value = 15659200
$ws.Range("J48").Value = 15382800
$ws.Range("D49").Value = 82100
$ws.Range("E49").Value = 64900
$ws.Range("F49").Value = 36700
$ws.Range("H49").Value = 362600
$ws.Range("I49").Value = 71200
$ws.Range("J49").Value = 74200
$ws.Range("D52").Value = 820400
$ws.Range("E52").Value = 746100
$ws.Range("F52").Value = 676700
$ws.Range("G52").Value = 628500
$ws.Range("H52").Value = 729600
$ws.Range("I52").Value = 928900
$ws.Range("J52").Value = 504300
$ws.Range("D54").Value = 32459900
$ws.Range("E54").Value = 29747600
$ws.Range("F54").Value = 27602600
$ws.Range("G54").Value = 28151600
$ws.Range("H54").Value = 24509200
$ws.Range("I54").Value = 21141600
$ws.Range("J54").Value = 19206000
$ws.Range("D57").Value = 315400
$ws.Range("E57").Value = 282400
$ws.Range("F57").Value = 371000
$ws.Range("G57").Value = 245900
$ws.Range("H57").Value = 1637700
$ws.Range("I57").Value = 270800
$ws.Range("J57").Value = 422500
$ws.Range("D58").Value = 5329300
$ws.Range("E58").Value = 5259800
$ws.Range("F58").Value = 5404800
$ws.Range("G58").Value = 4002800
$ws.Range("H58").Value = 2678900
$ws.Range("I58").Value = 7240300
$ws.Range("J58").Value = 3053200
$ws.Range("D59").Value = 4681300
$ws.Range("E59").Value = 4539600
$ws.Range("F59").Value = 3950200
$ws.Range("G59").Value = 3778200
$ws.Range("H59").Value = 6223300
$ws.Range("I59").Value = 3935400
$ws.Range("J59").Value = 3115600
$ws.Range("D60").Value = 10325900
$ws.Range("E60").Value = 10081800
$ws.Range("F60").Value = 9726000
$ws.Range("G60").Value = 8026900
$ws.Range("H60").Value = 7303700
$ws.Range("I60").Value = 7232200
$ws.Range("J60").Value = 6591300
$ws.Range("D61").Value = 11917600
$ws.Range("E61").Value = 10727800
$ws.Range("F61").Value = 9690000
$ws.Range("G61").Value = 12761000
$ws.Range("H61").Value = 5527700
$ws.Range("I61").Value = 7356200
$ws.Range("J61").Value = 6395000
$ws.Range("D62").Value = 934400
$ws.Range("E62").Value = 779000
$ws.Range("F62").Value = 821900
$ws.Range("G62").Value = 760500
$ws.Range("H62").Value = 6013300
$ws.Range("I62").Value = 940200
$ws.Range("J62").Value = 613200
$ws.Range("D66").Value = 25048900
$ws.Range("E66").Value = 23298300
$ws.Range("F66").Value = 21808000
$ws.Range("G66").Value = 22846200
$ws.Range("H66").Value = 19442600
$ws.Range("I66").Value = 16268000
$ws.Range("J66").Value = 14430900
$ws.Range("D72").Value = 3619900
$ws.Range("E72").Value = 2864000
$ws.Range("F72").Value = 2232200
$ws.Range("G72").Value = 1744600
$ws.Range("H72").Value = 3008700
$ws.Range("I72").Value = 1497800
$ws.Range("J72").Value = 1217700
$ws.Range("D76").Value = 7411000
$ws.Range("E76").Value = 6449300
$ws.Range("F76").Value = 5794700
$ws.Range("G76").Value = 5305400
$ws.Range("H76").Value = 5066600
$ws.Range("I76").Value = 4873600
$ws.Range("J76").Value = 4775100
$ws.Range("D81").Value = 884700
$ws.Range("E81").Value = 748600
$ws.Range("F81").Value = 554500
$ws.Range("G81").Value = 263700
$ws.Range("H81").Value = 294700
$ws.Range("I81").Value = 388700
$ws.Range("J81").Value = 758400
$ws.Range("D83").Value = 1977300
$ws.Range("E83").Value = 1892200
$ws.Range("F83").Value = 1779900
$ws.Range("G83").Value = 1630100
$ws.Range("H83").Value = 1408900
$ws.Range("I83").Value = 1237400
$ws.Range("J83").Value = 1152000
$ws.Range("D89").Value = 2631600
$ws.Range("E89").Value = 3526800
$ws.Range("F89").Value = 3522400
$ws.Range("G89").Value = 2013900
$ws.Range("H89").Value = 1440000
$ws.Range("I89").Value = 1737000
$ws.Range("J89").Value = 1863600
$ws.Range("D91").Value = -2054900
$ws.Range("E91").Value = -2814900
$ws.Range("F91").Value = -1801500
$ws.Range("G91").Value = -1283600
$ws.Range("H91").Value = -1826600
$ws.Range("I91").Value = -2334900
$ws.Range("J91").Value = -2973800
$ws.Range("D94").Value = -1222300
$ws.Range("E94").Value = -2337500
$ws.Range("F94").Value = -1028600
$ws.Range("G94").Value = -1448500
$ws.Range("H94").Value = -1811300
$ws.Range("I94").Value = -1803600
$ws.Range("J94").Value = -3258600
$ws.Range("D96").Value = -145700
$ws.Range("E96").Value = -116500
$ws.Range("F96").Value = -58300
$ws.Range("G96").Value = -58300
$ws.Range("H96").Value = -72900
$ws.Range("I96").Value = -291500
$ws.Range("D100").Value = -1008600
$ws.Range("E100").Value = -1255400
$ws.Range("F100").Value = -4110200
$ws.Range("G100").Value = -19400
$ws.Range("H100").Value = 618600
$ws.Range("I100").Value = 99100
$ws.Range("J100").Value = 1314800
$ws.Range("D101").Value = -3900
$ws.Range("E101").Value = 5500
$ws.Range("F101").Value = 5600
$ws.Range("D102").Value = 396800
$ws.Range("E102").Value = -60600
$ws.Range("F102").Value = -1610800
$ws.Range("G102").Value = 544100
$ws.Range("H102").Value = 247300
$ws.Range("I102").Value = 32500
$ws.Range("J102").Value = -80300
